$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D1").Value = "Locator Type"

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 4).Value = "CSS"
}

$ws.Range("D2:D21").Select()
